$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.243.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.13%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.365.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.41%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "552.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.34%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.80%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.106"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.93%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.56%  "

# Row 11
$ws.Range("E11").Value = "  -1.31%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.355"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.84%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.67%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.784.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.38%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.124.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.04%  "

# Row 16
$ws.Range("E16").Value = "  +1.84%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.424.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.80%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.87%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.76%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "331.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.28%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.92%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "63.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.53%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.169"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.45%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.16%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.51%  "

# Row 27
$ws.Range("E27").Value = "  -6.46%  "

# Row 28
$ws.Range("E28").Value = "  -0.23%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.05%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0742"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.44%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.32%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.67%  "

# Row 33
$ws.Range("E33").Value = "  -3.70%  "

# Row 34
$ws.Range("E34").Value = "  -0.01%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.03%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.09%  "

# Row 37
$ws.Range("E37").Value = "  -1.67%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.438"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +16.47%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "40.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.12%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.94%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.68%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "141.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.72%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "288.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.23%  "

# Row 44
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0960"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.40%  "

# Row 45
$ws.Range("B45").Value = "Polygon"
$ws.Range("C45").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.418"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.50%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0516"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.51%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.568"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.32%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.51%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0223"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.19%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.32%  "

# Row 51
$ws.Range("E51").Value = "  +0.13%  "
